# Restore C10 on the "Rules" worksheet from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
